$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the added "Outliers_MAD" columns
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Match the header formatting used by the existing header row (bold,
# centered, bordered) by copying the format from an existing header cell.
$ws.Range("C1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# Populate the new boolean "Outliers_MAD" columns for all data rows (2-12)
# with FALSE, matching the diff.
$ws.Range("F2:H12").Value = $false

$wb.Save()
